# Update "想去人数" (F) and "最低票价" (G) figures for a handful of events
# on both the "展览" and "全部类型" worksheets, matching output regenerated
# at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row -> (F new value, G new value or $null to leave unchanged) for the
# "展览" sheet.
$sheetExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @(
    @{ Row = 2;  F = 1338;  G = 70 },
    @{ Row = 4;  F = 14650 },
    @{ Row = 5;  F = 17876 },
    @{ Row = 7;  F = 75 },
    @{ Row = 17; F = 158 },
    @{ Row = 19; F = 1346 },
    @{ Row = 24; F = 7375 },
    @{ Row = 27; F = 43 },
    @{ Row = 28; F = 1185 },
    @{ Row = 30; F = 5880 },
    @{ Row = 33; F = 146 },
    @{ Row = 35; F = 232 },
    @{ Row = 36; F = 5132 }
)

foreach ($u in $expoUpdates) {
    $sheetExpo.Cells.Item($u.Row, 6).Value = $u.F
    if ($u.ContainsKey("G")) {
        $sheetExpo.Cells.Item($u.Row, 7).Value = $u.G
    }
}

# Same events, different row offsets on the "全部类型" sheet.
$sheetAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @(
    @{ Row = 2;  F = 1338;  G = 70 },
    @{ Row = 4;  F = 14650 },
    @{ Row = 5;  F = 17876 },
    @{ Row = 7;  F = 75 },
    @{ Row = 17; F = 158 },
    @{ Row = 19; F = 1346 },
    @{ Row = 25; F = 7375 },
    @{ Row = 28; F = 43 },
    @{ Row = 29; F = 1185 },
    @{ Row = 32; F = 5880 },
    @{ Row = 35; F = 146 },
    @{ Row = 37; F = 232 },
    @{ Row = 38; F = 5132 }
)

foreach ($u in $allUpdates) {
    $sheetAll.Cells.Item($u.Row, 6).Value = $u.F
    if ($u.ContainsKey("G")) {
        $sheetAll.Cells.Item($u.Row, 7).Value = $u.G
    }
}
